$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Survey Data")
$ws.Rows(11).Insert()
$ws.Rows(21).Insert()
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$ser1 = $chart.SeriesCollection().Item(1)
Write-Host "ser1 Formula after inserts:" $ser1.Formula
